# Auto-generated Excel COM-interop edit script
# Applies the scheduled-runner data refresh described in the diff:
# per-row currentAveragePrice / LevePrice / LeveProfit figures across the 8 crafting-class sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 324.45456
$ws.Range("I33").Value = 324.45456
$ws.Range("K33").Value = 324.45456
$ws.Range("M33").Value = -95.45456000000001
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H138").Value = 2643.8818
$ws.Range("I138").Value = 1264.619
$ws.Range("J138").Value = 3046.1667
$ws.Range("K138").Value = 3793.857
$ws.Range("L138").Value = 9138.500100000001
$ws.Range("M138").Value = 1346.143
$ws.Range("N138").Value = -19418.5001
$ws.Range("H141").Value = 16120.429
$ws.Range("I141").Value = 18157.166
$ws.Range("J141").Value = 3900
$ws.Range("K141").Value = 54471.49800000001
$ws.Range("L141").Value = 11700
$ws.Range("M141").Value = -49291.49800000001
$ws.Range("N141").Value = -22060

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 28810.666
$ws.Range("J76").Value = 28810.666
$ws.Range("L76").Value = 28810.666
$ws.Range("N76").Value = -29486.666
$ws.Range("H79").Value = 28810.666
$ws.Range("J79").Value = 28810.666
$ws.Range("L79").Value = 28810.666
$ws.Range("N79").Value = -31150.666
$ws.Range("H137").Value = 39035
$ws.Range("J137").Value = 40882
$ws.Range("L137").Value = 40882
$ws.Range("N137").Value = -51082

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 48000
$ws.Range("J59").Value = 48000
$ws.Range("L59").Value = 48000
$ws.Range("N59").Value = -49694
$ws.Range("H137").Value = 45550
$ws.Range("J137").Value = 45550
$ws.Range("L137").Value = 45550
$ws.Range("N137").Value = -55750

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H88").Value = 36228.4
$ws.Range("J88").Value = 36228.4
$ws.Range("L88").Value = 36228.4
$ws.Range("N88").Value = -37040.4
$ws.Range("H91").Value = 36228.4
$ws.Range("J91").Value = 36228.4
$ws.Range("L91").Value = 36228.4
$ws.Range("N91").Value = -39036.4
$ws.Range("H105").Value = 1707.6666
$ws.Range("I105").Value = 1300.4286
$ws.Range("J105").Value = 2277.8
$ws.Range("K105").Value = 1300.4286
$ws.Range("L105").Value = 2277.8
$ws.Range("M105").Value = 446.5714
$ws.Range("N105").Value = -5771.8
$ws.Range("H134").Value = 7457.421
$ws.Range("I134").Value = 9664.75
$ws.Range("J134").Value = 3673.4285
$ws.Range("K134").Value = 28994.25
$ws.Range("L134").Value = 11020.2855
$ws.Range("M134").Value = -26459.25
$ws.Range("N134").Value = -16090.2855

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 22288.8
$ws.Range("J34").Value = 12543.111
$ws.Range("L34").Value = 37629.333
$ws.Range("N34").Value = -37797.333
$ws.Range("H39").Value = 29333.334
$ws.Range("J39").Value = 40000
$ws.Range("L39").Value = 120000
$ws.Range("N39").Value = -120588
$ws.Range("H55").Value = 5500
$ws.Range("J55").Value = 5500
$ws.Range("L55").Value = 16500
$ws.Range("N55").Value = -16854
$ws.Range("H60").Value = 14644
$ws.Range("I60").Value = 350.6
$ws.Range("J60").Value = 38466.332
$ws.Range("K60").Value = 1051.8
$ws.Range("L60").Value = 115398.996
$ws.Range("M60").Value = -800.8000000000002
$ws.Range("N60").Value = -115900.996
$ws.Range("H68").Value = 1203
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 1203
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H92").Value = 781.6667
$ws.Range("J92").Value = 700
$ws.Range("L92").Value = 2100
$ws.Range("N92").Value = -4596

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10863
$ws.Range("I5").Value = 10863
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 10863
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -10751
$ws.Range("N5").ClearContents()
$ws.Range("H46").Value = 32626
$ws.Range("J46").Value = 35151.2
$ws.Range("L46").Value = 35151.2
$ws.Range("N46").Value = -35463.2
$ws.Range("H97").Value = 2249.5
$ws.Range("I97").Value = 1999.3334
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 1999.3334
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -1503.3334
$ws.Range("N97").Value = -3992
$ws.Range("H137").Value = 40530
$ws.Range("J137").Value = 40530
$ws.Range("L137").Value = 40530
$ws.Range("N137").Value = -50730

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2537.1428
$ws.Range("I2").Value = 750
$ws.Range("J2").Value = 2835
$ws.Range("K2").Value = 750
$ws.Range("L2").Value = 2835
$ws.Range("M2").Value = -638
$ws.Range("N2").Value = -3059
$ws.Range("H46").Value = 1948.2
$ws.Range("I46").Value = 1613.2667
$ws.Range("J46").Value = 2450.6
$ws.Range("K46").Value = 1613.2667
$ws.Range("L46").Value = 2450.6
$ws.Range("M46").Value = -1425.2667
$ws.Range("N46").Value = -2826.6
$ws.Range("H47").Value = 37499.5
$ws.Range("J47").Value = 37499.5
$ws.Range("L47").Value = 37499.5
$ws.Range("N47").Value = -38479.5
$ws.Range("H52").Value = 37499.5
$ws.Range("J52").Value = 37499.5
$ws.Range("L52").Value = 37499.5
$ws.Range("N52").Value = -37965.5
$ws.Range("H68").Value = 693.79
$ws.Range("I68").Value = 695.71716
$ws.Range("J68").Value = 503
$ws.Range("K68").Value = 695.71716
$ws.Range("L68").Value = 503
$ws.Range("M68").Value = 53.28283999999996
$ws.Range("N68").Value = -2001
$ws.Range("H71").Value = 693.79
$ws.Range("I71").Value = 695.71716
$ws.Range("J71").Value = 503
$ws.Range("K71").Value = 3478.5858
$ws.Range("L71").Value = 2515
$ws.Range("M71").Value = 265.4141999999997
$ws.Range("N71").Value = -10003

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 76910
$ws.Range("J46").Value = 76910
$ws.Range("L46").Value = 76910
$ws.Range("N46").Value = -77372
$ws.Range("H134").Value = 76910
$ws.Range("J134").Value = 76910
$ws.Range("L134").Value = 230730
$ws.Range("N134").Value = -235800
$ws.Range("H136").Value = 15786.917
$ws.Range("I136").Value = 27116.25
$ws.Range("K136").Value = 81348.75
$ws.Range("M136").Value = -78798.75
